$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F3").Value = 58
$ws1.Range("F4").Value = 78
$ws1.Range("F5").Value = 36
$ws1.Range("F6").Value = 36
$ws1.Range("F7").Value = 2640
$ws1.Range("F8").Value = 1152
$ws1.Range("F9").Value = 236
$ws1.Range("F10").Value = 95
$ws1.Range("F11").Value = 5950
$ws1.Range("F13").Value = 235
$ws1.Range("F14").Value = 586
$ws1.Range("F15").Value = 11611
$ws1.Range("F16").Value = 11857
$ws1.Range("F18").Value = 80
$ws1.Range("F21").Value = 63
$ws1.Range("F22").Value = 30

$ws4.Range("F3").Value = 58
$ws4.Range("F4").Value = 78
$ws4.Range("F5").Value = 36
$ws4.Range("F6").Value = 36
$ws4.Range("F7").Value = 2640
$ws4.Range("F9").Value = 1152
$ws4.Range("F10").Value = 236
$ws4.Range("F11").Value = 95
$ws4.Range("F12").Value = 5950
$ws4.Range("F14").Value = 235
$ws4.Range("F15").Value = 586
$ws4.Range("F16").Value = 11611
$ws4.Range("F17").Value = 11857
$ws4.Range("F19").Value = 80
$ws4.Range("F22").Value = 63
$ws4.Range("F23").Value = 30
